# Normalize the District column (G) values for the remaining rows that
# still contain school/taluk names instead of the cleaned district name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 9, 18, 20, 22, 28, 29, 34, 36, 38, 39, 42, 47, 48, 50)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Kalaburagi (Gulbarga)"
}
